$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 781.5
$ws.Range("I28").Value = 466.09525
$ws.Range("J28").Value = 1727.7142
$ws.Range("K28").Value = 466.09525
$ws.Range("L28").Value = 1727.7142
$ws.Range("M28").Value = 18.90474999999998
$ws.Range("N28").Value = -2697.7142
$ws.Range("H32").Value = 3362.8857
$ws.Range("J32").Value = 3575
$ws.Range("L32").Value = 3575
$ws.Range("N32").Value = -4227
$ws.Range("H62").Value = 7561.952
$ws.Range("I62").Value = 7618.727
$ws.Range("K62").Value = 7618.727
$ws.Range("M62").Value = -6994.727
$ws.Range("H65").Value = 7561.952
$ws.Range("I65").Value = 7618.727
$ws.Range("K65").Value = 38093.635
$ws.Range("M65").Value = -34973.635
$ws.Range("H86").Value = 10207.917
$ws.Range("I86").Value = 14249.167
$ws.Range("K86").Value = 14249.167
$ws.Range("M86").Value = -13126.167
$ws.Range("H89").Value = 10207.917
$ws.Range("I89").Value = 14249.167
$ws.Range("K89").Value = 71245.83499999999
$ws.Range("M89").Value = -65629.83499999999
$ws.Range("H107").Value = 47762012
$ws.Range("I107").Value = 83333460
$ws.Range("J107").Value = 333424.34
$ws.Range("K107").Value = 83333460
$ws.Range("L107").Value = 333424.34
$ws.Range("M107").Value = -83331540
$ws.Range("N107").Value = -337264.34
$ws.Range("H111").Value = 6947334
$ws.Range("I111").Value = 11113347
$ws.Range("K111").Value = 33340041
$ws.Range("M111").Value = -33336974
$ws.Range("H129").Value = 21741632
$ws.Range("I129").Value = 71429780
$ws.Range("K129").Value = 214289340
$ws.Range("M129").Value = -214284340
$ws.Range("H132").Value = 2351.04
$ws.Range("I132").Value = 2365.239
$ws.Range("K132").Value = 7095.717000000001
$ws.Range("M132").Value = -4565.717000000001
$ws.Range("H136").Value = 57000
$ws.Range("J136").Value = 57000
$ws.Range("L136").Value = 57000
$ws.Range("N136").Value = -67200
$ws.Range("H137").Value = 47997.51
$ws.Range("I137").Value = 59304.29
$ws.Range("K137").Value = 177912.87
$ws.Range("M137").Value = -175362.87
$ws.Range("H138").Value = 4276.524
$ws.Range("I138").Value = 3364
$ws.Range("J138").Value = 4732.7856
$ws.Range("K138").Value = 10092
$ws.Range("L138").Value = 14198.3568
$ws.Range("M138").Value = -4952
$ws.Range("N138").Value = -24478.3568
$ws.Range("H141").Value = 3183.6667
$ws.Range("I141").Value = 3172.7
$ws.Range("K141").Value = 9518.099999999999
$ws.Range("M141").Value = -4338.099999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 2421
$ws.Range("I31").Value = 2421
$ws.Range("K31").Value = 2421
$ws.Range("M31").Value = -2127
$ws.Range("H32").Value = 7350.804
$ws.Range("I32").Value = 6073.9766
$ws.Range("J32").Value = 22855.143
$ws.Range("K32").Value = 6073.9766
$ws.Range("L32").Value = 22855.143
$ws.Range("M32").Value = -5786.9766
$ws.Range("N32").Value = -23429.143
$ws.Range("H45").Value = 14290815
$ws.Range("I45").Value = 35715784
$ws.Range("K45").Value = 35715784
$ws.Range("M45").Value = -35715407
$ws.Range("H74").Value = 86869.5
$ws.Range("I74").Value = 14838.444
$ws.Range("J74").Value = 302962.66
$ws.Range("K74").Value = 14838.444
$ws.Range("L74").Value = 302962.66
$ws.Range("M74").Value = -13964.444
$ws.Range("N74").Value = -304710.66
$ws.Range("H77").Value = 86869.5
$ws.Range("I77").Value = 14838.444
$ws.Range("J77").Value = 302962.66
$ws.Range("K77").Value = 74192.22
$ws.Range("L77").Value = 1514813.3
$ws.Range("M77").Value = -69824.22
$ws.Range("N77").Value = -1523549.3
$ws.Range("H132").Value = 27882.076
$ws.Range("I132").Value = 5472.32
$ws.Range("K132").Value = 16416.96
$ws.Range("M132").Value = -13886.96
$ws.Range("H138").Value = 109998.5
$ws.Range("J138").Value = 109998.5
$ws.Range("L138").Value = 109998.5
$ws.Range("N138").Value = -120278.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 17750
$ws.Range("I102").Value = 5500
$ws.Range("K102").Value = 5500
$ws.Range("M102").Value = -2255
$ws.Range("H119").Value = 83196.336
$ws.Range("J119").Value = 83196.336
$ws.Range("L119").Value = 83196.336
$ws.Range("N119").Value = -92872.336
$ws.Range("H134").Value = 10269.883
$ws.Range("I134").Value = 7780
$ws.Range("J134").Value = 35998.668
$ws.Range("K134").Value = 23340
$ws.Range("L134").Value = 107996.004
$ws.Range("M134").Value = -20805
$ws.Range("N134").Value = -113066.004

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 959.3
$ws.Range("I22").Value = 1148.8667
$ws.Range("J22").Value = 390.6
$ws.Range("K22").Value = 1148.8667
$ws.Range("L22").Value = 390.6
$ws.Range("M22").Value = -798.8667
$ws.Range("N22").Value = -1090.6
$ws.Range("H31").Value = 28722.62
$ws.Range("I31").Value = 11207.615
$ws.Range("K31").Value = 11207.615
$ws.Range("M31").Value = -10912.615
$ws.Range("H34").Value = 28722.62
$ws.Range("I34").Value = 11207.615
$ws.Range("K34").Value = 11207.615
$ws.Range("M34").Value = -11005.615
$ws.Range("H58").Value = 4506.575
$ws.Range("I58").Value = 5352.32
$ws.Range("J58").Value = 3097
$ws.Range("K58").Value = 5352.32
$ws.Range("L58").Value = 3097
$ws.Range("M58").Value = -5149.32
$ws.Range("N58").Value = -3503
$ws.Range("H136").Value = 4506.575
$ws.Range("I136").Value = 5352.32
$ws.Range("J136").Value = 3097
$ws.Range("K136").Value = 16056.96
$ws.Range("L136").Value = 9291
$ws.Range("M136").Value = -13506.96
$ws.Range("N136").Value = -14391

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4065.6191
$ws.Range("J113").Value = 2231.2222
$ws.Range("L113").Value = 6693.6666
$ws.Range("N113").Value = -11033.6666
$ws.Range("H128").Value = 88996
$ws.Range("I128").Value = 88996
$ws.Range("K128").Value = 266988
$ws.Range("M128").Value = -262008
$ws.Range("H133").Value = 4681.2856
$ws.Range("I133").Value = 3889.5
$ws.Range("J133").Value = 4998
$ws.Range("K133").Value = 11668.5
$ws.Range("L133").Value = 14994
$ws.Range("M133").Value = -6608.5
$ws.Range("N133").Value = -25114
$ws.Range("H141").Value = 3704.2856
$ws.Range("I141").Value = 3488.3333
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 10464.9999
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -5284.999899999999
$ws.Range("N141").Value = -25360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 8057.149
$ws.Range("I132").Value = 6385.3423
$ws.Range("K132").Value = 19156.0269
$ws.Range("M132").Value = -16626.0269

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5263.9805
$ws.Range("I7").Value = 4392.448
$ws.Range("J7").Value = 6412.8184
$ws.Range("K7").Value = 4392.448
$ws.Range("L7").Value = 6412.8184
$ws.Range("M7").Value = -4280.448
$ws.Range("N7").Value = -6636.8184
$ws.Range("H126").Value = 5263.9805
$ws.Range("I126").Value = 4392.448
$ws.Range("J126").Value = 6412.8184
$ws.Range("K126").Value = 13177.344
$ws.Range("L126").Value = 19238.4552
$ws.Range("M126").Value = -10707.344
$ws.Range("N126").Value = -24178.4552
$ws.Range("H136").Value = 78745.484
$ws.Range("I136").Value = 127384.69
$ws.Range("K136").Value = 382154.07
$ws.Range("M136").Value = -379604.07

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 13350.4
$ws.Range("J74").Value = 13350.4
$ws.Range("L74").Value = 13350.4
$ws.Range("N74").Value = -15222.4
$ws.Range("H77").Value = 13350.4
$ws.Range("J77").Value = 13350.4
$ws.Range("L77").Value = 40051.2
$ws.Range("N77").Value = -49411.2
$ws.Range("H132").Value = 15802949
$ws.Range("I132").Value = 18188914
$ws.Range("J132").Value = 1222056.1
$ws.Range("K132").Value = 54566742
$ws.Range("L132").Value = 3666168.3
$ws.Range("M132").Value = -54564212
$ws.Range("N132").Value = -3671228.3
$ws.Range("H136").Value = 4041.0164
$ws.Range("I136").Value = 4134.96
$ws.Range("K136").Value = 12404.88
$ws.Range("M136").Value = -9854.880000000001
